# Tracking Table.xlsx edit script
# Adds "Leyland Steam Van" (Heavy Goods) and five light-goods vans
# (Ford Transit Mk1, Ford Transit Mk3, Bedford Rascal, Ford Transit Custom,
# Mercedes-Benz Sprinter, Ford Thames) to the vehicle tracking table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: append "Ford Transit Mk1" as new row 35 (bottom of current table) ---
$ws.Cells.Item(35,1).Value = "Ford Transit Mk1"
$ws.Cells.Item(35,2).Value = 1965
$ws.Cells.Item(35,3).Value = 1
$ws.Cells.Item(35,4).Value = "Light Goods"
$ws.Cells.Item(35,5).Formula = "=IF(B35 > 1900, ((B35-1900)*10)+400+C35, ((B35-1730)*2)+C35)+VLOOKUP(D35,'ID Scheme'!`$A`$2:`$B`$4,2)"
$ws.Cells.Item(35,6).Value = 65
$ws.Cells.Item(35,7).Value = 8
$ws.Cells.Item(35,8).Formula = "=SQRT(F35*G35)/`$B`$1"
$ws.Cells.Item(35,8).NumberFormat = "0"
$ws.Cells.Item(35,9).Formula = "=H35*0.9"
$ws.Cells.Item(35,9).NumberFormat = "0"
$ws.Cells.Item(35,10).Value = "x"
$ws.Cells.Item(35,10).NumberFormat = "0"

# --- Step 2: append "Ford Transit Mk3" as new row 36 ---
$ws.Cells.Item(36,1).Value = "Ford Transit Mk3"
$ws.Cells.Item(36,2).Value = 1986
$ws.Cells.Item(36,3).Value = 1
$ws.Cells.Item(36,4).Value = "Light Goods"
$ws.Cells.Item(36,5).Formula = "=IF(B36 > 1900, ((B36-1900)*10)+400+C36, ((B36-1730)*2)+C36)+VLOOKUP(D36,'ID Scheme'!`$A`$2:`$B`$4,2)"
$ws.Cells.Item(36,6).Value = 80
$ws.Cells.Item(36,7).Value = 10
$ws.Cells.Item(36,8).Formula = "=SQRT(F36*G36)/`$B`$1"
$ws.Cells.Item(36,8).NumberFormat = "0"
$ws.Cells.Item(36,9).Formula = "=H36*0.9"
$ws.Cells.Item(36,9).NumberFormat = "0"
$ws.Cells.Item(36,10).Value = "x"
$ws.Cells.Item(36,10).NumberFormat = "0"

# --- Step 3: insert a new row at 23 for "Leyland Steam Van" (Heavy Goods) ---
# This pushes the existing rows 23-36 down to 24-37.
$ws.Rows.Item(23).Insert()

$ws.Cells.Item(23,1).Value = "Leyland Steam Van"
$ws.Cells.Item(23,2).Value = 1897
$ws.Cells.Item(23,3).Value = 1
$ws.Cells.Item(23,4).Value = "Heavy Goods"
$ws.Cells.Item(23,5).Formula = "=IF(B23 > 1900, ((B23-1900)*10)+400+C23, ((B23-1730)*2)+C23)+VLOOKUP(D23,'ID Scheme'!`$A`$2:`$B`$4,2)"
$ws.Cells.Item(23,6).Value = 18
$ws.Cells.Item(23,7).Value = 14
$ws.Cells.Item(23,8).Formula = "=SQRT(F23*G23)/`$B`$1"
$ws.Cells.Item(23,8).NumberFormat = "0"
$ws.Cells.Item(23,9).Formula = "=H23*0.9"
$ws.Cells.Item(23,9).NumberFormat = "0"
$ws.Cells.Item(23,10).Value = "x"
$ws.Cells.Item(23,10).NumberFormat = "0"

# --- Step 4: append "Bedford Rascal" as new row 38 ---
$ws.Cells.Item(38,1).Value = "Bedford Rascal"
$ws.Cells.Item(38,2).Value = 1986
$ws.Cells.Item(38,3).Value = 2
$ws.Cells.Item(38,4).Value = "Light Goods"
$ws.Cells.Item(38,5).Formula = "=IF(B38 > 1900, ((B38-1900)*10)+400+C38, ((B38-1730)*2)+C38)+VLOOKUP(D38,'ID Scheme'!`$A`$2:`$B`$4,2)"
$ws.Cells.Item(38,6).Value = 65
$ws.Cells.Item(38,7).Value = 4
$ws.Cells.Item(38,8).Formula = "=SQRT(F38*G38)/`$B`$1"
$ws.Cells.Item(38,8).NumberFormat = "0"
$ws.Cells.Item(38,9).Formula = "=H38*0.9"
$ws.Cells.Item(38,9).NumberFormat = "0"
$ws.Cells.Item(38,10).Value = "x"
$ws.Cells.Item(38,10).NumberFormat = "0"

# --- Step 5: append "Ford Transit Custom" as new row 39 ---
$ws.Cells.Item(39,1).Value = "Ford Transit Custom"
$ws.Cells.Item(39,2).Value = 2012
$ws.Cells.Item(39,3).Value = 1
$ws.Cells.Item(39,4).Value = "Light Goods"
$ws.Cells.Item(39,5).Formula = "=IF(B39 > 1900, ((B39-1900)*10)+400+C39, ((B39-1730)*2)+C39)+VLOOKUP(D39,'ID Scheme'!`$A`$2:`$B`$4,2)"
$ws.Cells.Item(39,6).Value = 92
$ws.Cells.Item(39,7).Value = 10
$ws.Cells.Item(39,8).Formula = "=SQRT(F39*G39)/`$B`$1"
$ws.Cells.Item(39,8).NumberFormat = "0"
$ws.Cells.Item(39,9).Formula = "=H39*0.9"
$ws.Cells.Item(39,9).NumberFormat = "0"
$ws.Cells.Item(39,10).Value = "x"
$ws.Cells.Item(39,10).NumberFormat = "0"

# --- Step 6: append "Mercedes-Benz Sprinter" as new row 40 ---
$ws.Cells.Item(40,1).Value = "Mercedes-Benz Sprinter"
$ws.Cells.Item(40,2).Value = 2006
$ws.Cells.Item(40,3).Value = 1
$ws.Cells.Item(40,4).Value = "Light Goods"
$ws.Cells.Item(40,5).Formula = "=IF(B40 > 1900, ((B40-1900)*10)+400+C40, ((B40-1730)*2)+C40)+VLOOKUP(D40,'ID Scheme'!`$A`$2:`$B`$4,2)"
$ws.Cells.Item(40,6).Value = 85
$ws.Cells.Item(40,7).Value = 12
$ws.Cells.Item(40,8).Formula = "=SQRT(F40*G40)/`$B`$1"
$ws.Cells.Item(40,8).NumberFormat = "0"
$ws.Cells.Item(40,9).Formula = "=H40*0.9"
$ws.Cells.Item(40,9).NumberFormat = "0"
$ws.Cells.Item(40,10).Value = "x"
$ws.Cells.Item(40,10).NumberFormat = "0"

# --- Step 7: append "Ford Thames" as new row 41 ---
$ws.Cells.Item(41,1).Value = "Ford Thames"
$ws.Cells.Item(41,2).Value = 1957
$ws.Cells.Item(41,3).Value = 1
$ws.Cells.Item(41,4).Value = "Light Goods"
$ws.Cells.Item(41,5).Formula = "=IF(B41 > 1900, ((B41-1900)*10)+400+C41, ((B41-1730)*2)+C41)+VLOOKUP(D41,'ID Scheme'!`$A`$2:`$B`$4,2)"
$ws.Cells.Item(41,6).Value = 60
$ws.Cells.Item(41,7).Value = 6
$ws.Cells.Item(41,8).Formula = "=SQRT(F41*G41)/`$B`$1"
$ws.Cells.Item(41,8).NumberFormat = "0"
$ws.Cells.Item(41,9).Formula = "=H41*0.9"
$ws.Cells.Item(41,9).NumberFormat = "0"
$ws.Cells.Item(41,10).Value = "x"
$ws.Cells.Item(41,10).NumberFormat = "0"

# --- Update sheet view to match the final state ---
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
$ws.Range("H35").Select()
